$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Burndown Chart Template")
$ws.Range("E17").Value = 3
